# Update QA STATUS/COMMENT values per new STATUS tracking scheme
# Valid STATUS values going forward: ISSUE, NO ISSUE, BLOCKED

$wb = $excel.ActiveWorkbook

# --- Sheet1 ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("E3").Value = "ISSUE"
$ws1.Range("F3").Value = "Agree - typo"

$ws1.Range("E5").Value = "NO ISSUE"
$ws1.Range("F5").Value = "Good"

$ws1.Range("E6").Value = "ISSUE"

# --- Sheet2 ---
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("E2").Value = "NO ISSUE"

$ws2.Range("E3").Value = "NO ISSUE"
$ws2.Range("F3").Value = "Verified"

$ws2.Range("E4").Value = "NO ISSUE"
$ws2.Range("F4").Value = "Correct"
